# daily auto push: 2026-01-23 06:52 UTC
# Insert a new data row for 2026/01/23 (13:00, rank 10) into the "sei3" log
# sheet at row 682, pushing the existing rows (2026/12/29 onward) down by
# one. Dimension (A1:D723 -> A1:D724) is updated automatically by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything at/after row 682 down by one row.
$ws.Rows.Item(682).Insert()

# Populate the freshly inserted row. The date column holds a literal text
# string (e.g. "2026/01/23"), not a real date value, elsewhere in this
# sheet, so prefix with an apostrophe to force text entry and then reset
# the cell style back to Normal/General (Excel would otherwise stamp the
# cell with a "quote prefix" number format just because of the apostrophe).
$ws.Range("A682").Value = "'2026/01/23"
$ws.Range("A682").Style = "Normal"

$ws.Range("B682").Value = "金"
$ws.Range("C682").Value = 13
$ws.Range("D682").Value = 10
